$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need the column
# pre-formatted as Text so Excel stores the literal string (matching
# the source data, which is all text) instead of converting it to a
# numeric value.
$ws.Range("D5","D6","D7","D9","D10","D11","D12","D14","D21","D22","D23","D24","D27","D28","D31","D32","D33","D34","D37","D39","D40","D41","D42","D45","D47","D48","D50","D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.194.24'
$ws.Range("D3").Value = '2.283.03'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  -0.68%  '
$ws.Range("D5").Value = '113.78'
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").Value = '265.66'
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  -0.86%  '
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D9").Value = '0.612'
$ws.Range("E9").Value = '  -0.77%  '
$ws.Range("D10").Value = '47.75'
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").Value = '0.0926'
$ws.Range("E11").Value = '  -0.85%  '
$ws.Range("D12").Value = '9.09'
$ws.Range("E12").Value = '  +5.55%  '
$ws.Range("E13").Value = '  +1.50%  '
$ws.Range("D14").Value = '15.58'
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '2.625.10'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("D17").Value = '2.277.59'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").Value = '43.221.63'
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("E20").Value = '  +4.95%  '
$ws.Range("D21").Value = '71.62'
$ws.Range("E21").Value = '  -1.08%  '
$ws.Range("D22").Value = '2.53'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").Value = '232.45'
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '9.67'
$ws.Range("E24").Value = '  +1.71%  '
$ws.Range("E25").Value = '  +1.29%  '
$ws.Range("E26").Value = '  +1.33%  '
$ws.Range("D27").Value = '11.40'
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").Value = '40.83'
$ws.Range("E28").Value = '  -6.83%  '
$ws.Range("E29").Value = '  -2.30%  '
$ws.Range("E30").Value = '  -0.96%  '
$ws.Range("D31").Value = '172.61'
$ws.Range("E31").Value = '  -2.13%  '
$ws.Range("D32").Value = '21.43'
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("D33").Value = '0.0908'
$ws.Range("E33").Value = '  -2.00%  '
$ws.Range("D34").Value = '5.79'
$ws.Range("E34").Value = '  +5.36%  '
$ws.Range("E35").Value = '  +1.19%  '
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").Value = '3.96'
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("D39").Value = '0.104'
$ws.Range("E39").Value = '  -4.75%  '
$ws.Range("D40").Value = '2.68'
$ws.Range("E40").Value = '  +12.46%  '
$ws.Range("D41").Value = '77.90'
$ws.Range("E41").Value = '  +3.94%  '
$ws.Range("D42").Value = '14.00'
$ws.Range("E42").Value = '  +6.45%  '
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("E44").Value = '  +5.10%  '
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("E46").Value = '  -1.69%  '
$ws.Range("D47").Value = '8.72'
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("D48").Value = '104.46'
$ws.Range("E48").Value = '  +3.55%  '
$ws.Range("E49").Value = '  +1.75%  '
$ws.Range("D50").Value = '0.0998'
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("D51").Value = '0.439'
$ws.Range("E51").Value = '  -3.48%  '
